$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1
$ws.Range("D3").Value = 1

$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 2

$ws.Range("B7").Value = 4
$ws.Range("D7").Value = 2

$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 2

$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 4

$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 1
$ws.Range("E10").Value = 2
